# Daily attendance processing - 2026-01-19 19:14:50
# Normalizes the "Recorded By" (column G) entries so that "System" is
# listed after the other recorder(s) rather than first, for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2   = "system, backup@backdoor.com, System"
    3   = "dnasr281@gmail.com, System"
    4   = "backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    6   = "dnasr281@gmail.com, System"
    7   = "admin@admin.com, System"
    8   = "backup@backdoor.com, System"
    28  = "system, backup@backdoor.com, System"
    29  = "dnasr281@gmail.com, System"
    30  = "backup@backdoor.com, System"
    31  = "backup@backdoor.com, System"
    32  = "dnasr281@gmail.com, System"
    33  = "admin@admin.com, System"
    34  = "backup@backdoor.com, System"
    54  = "system, backup@backdoor.com, System"
    55  = "dnasr281@gmail.com, System"
    56  = "backup@backdoor.com, System"
    57  = "backup@backdoor.com, System"
    58  = "dnasr281@gmail.com, System"
    59  = "admin@admin.com, System"
    60  = "backup@backdoor.com, System"
    80  = "backup@backdoor.com, System"
    81  = "backup@backdoor.com, System"
    82  = "backup@backdoor.com, System"
    87  = "admin@admin.com, dnasr281@gmail.com"
    106 = "backup@backdoor.com, System"
    107 = "backup@backdoor.com, System"
    108 = "backup@backdoor.com, System"
    113 = "admin@admin.com, dnasr281@gmail.com"
    132 = "backup@backdoor.com, System"
    133 = "backup@backdoor.com, System"
    134 = "backup@backdoor.com, System"
    139 = "admin@admin.com, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
